$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("999")

# Row 3 (right-hand "Udhayam" block) is being cleared out - the 999 check
# for S.No 1 no longer has a matching tracked row on the right side.
# K3 keeps its date style/format but loses its value.
$ws.Range("I3").ClearContents()
$ws.Range("J3").ClearContents()
$ws.Range("K3").ClearContents()
$ws.Range("L3").ClearContents()

# New row 5: a third day of tracking data for both Kabeer (left block)
# and Udhayam (right block).
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "Kabeer"
$ws.Range("C5").NumberFormat = "m/d/yy"
$ws.Range("C5").Value = "3/16/2018"
$ws.Range("D5").Value = 10

$ws.Range("I5").Value = 3
$ws.Range("J5").Value = "Udhayam"
$ws.Range("K5").NumberFormat = "m/d/yy"
$ws.Range("K5").Value = "3/16/2018"
$ws.Range("L5").Value = 9

# Make the "999" sheet the active tab/selection, matching the saved
# view state in the workbook (it was Sheet1 before, now it's "999").
$ws.Activate()
$ws.Range("I3:L3").Select()
